# 2021 data cleaning complete
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the two new flag/description rows
$ws.Range("A6").Value = "m"
$ws.Range("B6").Value = "rain gague malfunction; values incorrect - recoded to na"
$ws.Range("A7").Value = "e"
$ws.Range("B7").Value = "data errant, recoded to na"

# Update selection to match the target (A8)
$ws.Range("A8").Select()
